$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.295968333333334
$ws.Range("H2").Value = 18.887905
$ws.Range("I2").Value = 0.5052862712055841
$ws.Range("J2").Value = 0.5052862712055841
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05828766666666666
$ws.Range("N2").Value = 0.174863
$ws.Range("O2").Value = 0.01080277125928955
$ws.Range("P2").Value = 0.01080277125928955
$ws.Range("Q2").Value = 0.3669773035572222
$ws.Range("R2").Value = 3.302795732015
$ws.Range("S2").Value = 0.00545849200829327
$ws.Range("T2").Value = 0.00545849200829327

$ws.Range("G3").Value = 6.295968333333334
$ws.Range("H3").Value = 18.887905
$ws.Range("I3").Value = 0.5052862712055841
$ws.Range("J3").Value = 0.5052862712055841
$ws.Range("M3").Value = 3.424957333333333
$ws.Range("O3").Value = 0.634766027887426
$ws.Range("P3").Value = 0.634766027887426
$ws.Range("Q3").Value = 21.56342291368444
$ws.Range("R3").Value = 194.07080622316
$ws.Range("S3").Value = 0.3207385593192173
$ws.Range("T3").Value = 0.3207385593192173

$ws.Range("G4").Value = 6.295968333333334
$ws.Range("H4").Value = 18.887905
$ws.Range("I4").Value = 0.5052862712055841
$ws.Range("J4").Value = 0.5052862712055841
$ws.Range("M4").Value = 1.912376666666667
$ws.Range("N4").Value = 5.737130000000001
$ws.Range("O4").Value = 0.3544312008532844
$ws.Range("P4").Value = 0.3544312008532844
$ws.Range("Q4").Value = 12.04026293473889
$ws.Range("R4").Value = 108.36236641265
$ws.Range("S4").Value = 0.1790892198780736
$ws.Range("T4").Value = 0.1790892198780736

$ws.Range("I5").Value = 0.2025983155648483
$ws.Range("J5").Value = 0.2025983155648483
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05828766666666666
$ws.Range("N5").Value = 0.174863
$ws.Range("O5").Value = 0.01080277125928955
$ws.Range("P5").Value = 0.01080277125928955
$ws.Range("Q5").Value = 0.1471422989067778
$ws.Range("R5").Value = 1.324280690161
$ws.Range("S5").Value = 0.002188623260564418
$ws.Range("T5").Value = 0.002188623260564419

$ws.Range("I6").Value = 0.2025983155648483
$ws.Range("J6").Value = 0.2025983155648483
$ws.Range("M6").Value = 3.424957333333333
$ws.Range("O6").Value = 0.634766027887426
$ws.Range("P6").Value = 0.634766027887426
$ws.Range("R6").Value = 77.814143549384
$ws.Range("S6").Value = 0.128602528027782
$ws.Range("T6").Value = 0.128602528027782

$ws.Range("I7").Value = 0.2025983155648483
$ws.Range("J7").Value = 0.2025983155648483
$ws.Range("M7").Value = 1.912376666666667
$ws.Range("N7").Value = 5.737130000000001
$ws.Range("O7").Value = 0.3544312008532844
$ws.Range("P7").Value = 0.3544312008532844
$ws.Range("Q7").Value = 4.827633617901111
$ws.Range("R7").Value = 43.44870256111
$ws.Range("S7").Value = 0.07180716427650186
$ws.Range("T7").Value = 0.07180716427650187

$ws.Range("G8").Value = 3.639816666666666
$ws.Range("H8").Value = 10.91945
$ws.Range("I8").Value = 0.2921154132295675
$ws.Range("J8").Value = 0.2921154132295676
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05828766666666666
$ws.Range("N8").Value = 0.174863
$ws.Range("O8").Value = 0.01080277125928955
$ws.Range("P8").Value = 0.01080277125928955
$ws.Range("Q8").Value = 0.2121564205944444
$ws.Range("R8").Value = 1.90940778535
$ws.Range("S8").Value = 0.003155655990431863
$ws.Range("T8").Value = 0.003155655990431864

$ws.Range("G9").Value = 3.639816666666666
$ws.Range("H9").Value = 10.91945
$ws.Range("I9").Value = 0.2921154132295675
$ws.Range("J9").Value = 0.2921154132295676
$ws.Range("M9").Value = 3.424957333333333
$ws.Range("O9").Value = 0.634766027887426
$ws.Range("P9").Value = 0.634766027887426
$ws.Range("Q9").Value = 12.46621678448889
$ws.Range("R9").Value = 112.1959510604
$ws.Range("S9").Value = 0.1854249405404266
$ws.Range("T9").Value = 0.1854249405404267

$ws.Range("G10").Value = 3.639816666666666
$ws.Range("H10").Value = 10.91945
$ws.Range("I10").Value = 0.2921154132295675
$ws.Range("J10").Value = 0.2921154132295676
$ws.Range("M10").Value = 1.912376666666667
$ws.Range("N10").Value = 5.737130000000001
$ws.Range("O10").Value = 0.3544312008532844
$ws.Range("P10").Value = 0.3544312008532844
$ws.Range("Q10").Value = 6.960700464277778
$ws.Range("R10").Value = 62.6463041785
$ws.Range("S10").Value = 0.103534816698709
$ws.Range("T10").Value = 0.1035348166987091

Write-Output "done"
